$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create the new sheet "3" as the last tab, as a copy-like rebuild of sheet "1" ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "3"

# Header row (reuses the same shared strings as sheet "1")
$ws3.Range("A1").Value = "Rok"
$ws3.Range("B1").Value = "Liczba Żubrów"
$ws3.Range("C1").Value = "Różnica "
$ws3.Range("D1").Value = "Procentowo"

# Year column
$ws3.Range("A2").Value = 2013
$ws3.Range("A3").Value = 2014
$ws3.Range("A4").Value = 2015
$ws3.Range("A5").Value = 2016
$ws3.Range("A6").Value = 2017
$ws3.Range("A7").Value = 2018
$ws3.Range("A8").Value = 2019
$ws3.Range("A9").Value = 2020
$ws3.Range("A10").Value = 2021
$ws3.Range("A11").Value = 2022

# "Liczba Żubrów" column - same formulas as sheet "1", except 2018 (B7) gains an
# extra "+230" term per the revised count.
$ws3.Range("B2").Formula = "=7+40+2+3+3+10+14+2+6+1+5+6+11+25+3+46+6+21+4+5+3+4+5+4+3+270+505+103+128"
$ws3.Range("B3").Formula = "=9+37+2+3+3+11+5+2+6+3+4+5+10+22+2+41+6+15+4+5+1+4+5+4+301+522+105+123+159"
$ws3.Range("B4").Formula = "=8+28+2+3+3+10+11+6+4+6+3+7+12+27+2+43+8+7+1+4+6+4+4+344+578+107+134+184"
$ws3.Range("B5").Formula = "=9+39+2+3+3+12+12+9+4+7+1+8+16+26+2+50+7+6+4+8+1+4+5+5+402+596+108+144+205"
$ws3.Range("B6").Formula = "=8+43+1+3+4+13+10+8+4+8+9+11+27+3+42+7+6+4+9+1+5+7+5+487+654+120+158+216"
$ws3.Range("B7").Formula = "=8+19+1+4+4+12+7+9+3+8+9+11+23+5+48+6+4+10+1+3+6+6+551+8+519+112+158+265+230"
$ws3.Range("B8").Formula = "=8+27+1+4+5+13+5+10+3+6+6+9+22+5+56+7+4+7+1+3+8+5+6+668+9+770+112+184+305"
$ws3.Range("B9").Formula = "=6+31+1+11+3+4+5+6+7+3+4+8+13+19+6+45+5+7+4+7+5+3+7+6+3+707+17+715+117+214+334"
$ws3.Range("B10").Formula = "=7+28+1+4+4+7+7+9+3+8+10+16+7+50+6+6+7+3+2+9+8+4+729+9+20+779+125+212+9+340"
$ws3.Range("B11").Value = 2603

# "Różnica " column: C3 stands alone, C4:C11 fill down as one shared formula
# (mirrors the structure already present on sheet "1").
$ws3.Range("C3").Formula = "=B3-B2"
$ws3.Range("C4:C11").Formula = "=B4-B3"

# Build the table over the populated range and rename it so the new sheet has
# its own table (avoids clashing with "Tabela3" on sheet "1").
$lo = $ws3.ListObjects.Add(1, $ws3.Range("A1:D11"), $null, 1)
$lo.Name = "Tabela32"

# "Procentowo" column - calculated from the new table's own name.
$ws3.Range("D2").Formula = "=(Tabela32[[#This Row],[Różnica ]]/Tabela32[[#This Row],[Liczba Żubrów]])*100%"
$ws3.Range("D3").Formula = "=(Tabela32[[#This Row],[Różnica ]]/Tabela32[[#This Row],[Liczba Żubrów]])*100%"
$ws3.Range("D4").Formula = "=(Tabela32[[#This Row],[Różnica ]]/Tabela32[[#This Row],[Liczba Żubrów]])*100%"
$ws3.Range("D5").Formula = "=(Tabela32[[#This Row],[Różnica ]]/Tabela32[[#This Row],[Liczba Żubrów]])*100%"
$ws3.Range("D6").Formula = "=(Tabela32[[#This Row],[Różnica ]]/Tabela32[[#This Row],[Liczba Żubrów]])*100%"
$ws3.Range("D7").Formula = "=(Tabela32[[#This Row],[Różnica ]]/Tabela32[[#This Row],[Liczba Żubrów]])*100%"
$ws3.Range("D8").Formula = "=(Tabela32[[#This Row],[Różnica ]]/Tabela32[[#This Row],[Liczba Żubrów]])*100%"
$ws3.Range("D9").Formula = "=(Tabela32[[#This Row],[Różnica ]]/Tabela32[[#This Row],[Liczba Żubrów]])*100%"
$ws3.Range("D10").Formula = "=(Tabela32[[#This Row],[Różnica ]]/Tabela32[[#This Row],[Liczba Żubrów]])*100%"
$ws3.Range("D11").Formula = "=(Tabela32[[#This Row],[Różnica ]]/Tabela32[[#This Row],[Liczba Żubrów]])*100%"

# Percentage number format for D3:D11, matching the "Procentowy" cell style
# already used on sheet "1".
$ws3.Range("D3:D11").NumberFormat = "0%"

# --- Selections / active tab ---
# Sheet "1" loses the tab focus and ends up with its whole table selected.
$ws1.Activate()
$ws1.Range("A1:D11").Select()

# Sheet "3" becomes the active tab, with B8 selected.
$ws3.Activate()
$ws3.Range("B8").Select()
